# "changed tarif to june"
# Sheet1 "БСК Тариф НП": the per-destination tariff table (rows 4-9) is
# re-priced and the destination groupings are re-split into two more rows
# (the table grows from 6 data rows to 8 data rows: rows 4-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Insert two new rows right after the current last data row (row 9),
#    copying formatting (borders/fonts/number formats) down from row 9
#    so the new rows look consistent with the existing table.
# ------------------------------------------------------------------
$ws.Rows("10:11").Insert()

$ws.Range("A9:Q9").Copy()
$ws.Range("A10:Q10").PasteSpecial(-4122)
$ws.Range("A9:Q9").Copy()
$ws.Range("A11:Q11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Extend the A4:A9 merged "Ростов-на-Дону" label down to A4:A11.
# ------------------------------------------------------------------
$ws.Range("A4:A9").UnMerge()
$ws.Range("A4:A11").Merge()
$ws.Range("A4").Value = "Ростов-на-Дону"

# ------------------------------------------------------------------
# 3. Re-write the destination labels (col B) and the June tariff
#    numbers (cols C:Q) for every row of the table, row by row.
# ------------------------------------------------------------------

# Row 4 - Донецк
$ws.Range("B4").Value = "Донецк "
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 600
$ws.Range("E4").Value = 34
$ws.Range("F4").Value = 20
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 19
$ws.Range("I4").Value = 18
$ws.Range("J4").Value = 18
$ws.Range("K4").Value = 600
$ws.Range("L4").Value = 8500
$ws.Range("M4").Value = 5000
$ws.Range("N4").Value = 5000
$ws.Range("O4").Value = 4750
$ws.Range("P4").Value = 4500
$ws.Range("Q4").Value = 4500

# Row 5 - Макеевка,
$ws.Range("B5").Value = "Макеевка, "
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 600
$ws.Range("E5").Value = 30
$ws.Range("F5").Value = 18
$ws.Range("G5").Value = 17
$ws.Range("H5").Value = 17
$ws.Range("I5").Value = 16
$ws.Range("J5").Value = 16
$ws.Range("K5").Value = 600
$ws.Range("L5").Value = 7500
$ws.Range("M5").Value = 4500
$ws.Range("N5").Value = 4250
$ws.Range("O5").Value = 4250
$ws.Range("P5").Value = 4000
$ws.Range("Q5").Value = 4000

# Row 6 - Горловка/ Енакиево/
$ws.Range("B6").Value = "Горловка/ Енакиево/ "
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 600
$ws.Range("E6").Value = 34
$ws.Range("F6").Value = 20
$ws.Range("G6").Value = 20
$ws.Range("H6").Value = 19
$ws.Range("I6").Value = 18
$ws.Range("J6").Value = 18
$ws.Range("K6").Value = 600
$ws.Range("L6").Value = 8500
$ws.Range("M6").Value = 5000
$ws.Range("N6").Value = 5000
$ws.Range("O6").Value = 4750
$ws.Range("P6").Value = 4500
$ws.Range("Q6").Value = 4500

# Row 7 - Шахтерск/ Снежное / Кр Луч
$ws.Range("B7").Value = "Шахтерск/ Снежное / Кр Луч"
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 600
$ws.Range("E7").Value = 34
$ws.Range("F7").Value = 20
$ws.Range("G7").Value = 20
$ws.Range("H7").Value = 19
$ws.Range("I7").Value = 17
$ws.Range("J7").Value = 16
$ws.Range("K7").Value = 600
$ws.Range("L7").Value = 8500
$ws.Range("M7").Value = 5000
$ws.Range("N7").Value = 5000
$ws.Range("O7").Value = 4750
$ws.Range("P7").Value = 4250
$ws.Range("Q7").Value = 4000

# Row 8 - Луганск
$ws.Range("B8").Value = "Луганск"
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 600
$ws.Range("E8").Value = 30
$ws.Range("F8").Value = 18
$ws.Range("G8").Value = 17
$ws.Range("H8").Value = 17
$ws.Range("I8").Value = 16
$ws.Range("J8").Value = 16
$ws.Range("K8").Value = 600
$ws.Range("L8").Value = 7500
$ws.Range("M8").Value = 4500
$ws.Range("N8").Value = 4250
$ws.Range("O8").Value = 4250
$ws.Range("P8").Value = 4000
$ws.Range("Q8").Value = 4000

# Row 9 - Стаханова/ Алчевск
$ws.Range("B9").Value = "Стаханова/ Алчевск"
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 600
$ws.Range("E9").Value = 34
$ws.Range("F9").Value = 20
$ws.Range("G9").Value = 20
$ws.Range("H9").Value = 19
$ws.Range("I9").Value = 18
$ws.Range("J9").Value = 18
$ws.Range("K9").Value = 600
$ws.Range("L9").Value = 8500
$ws.Range("M9").Value = 5000
$ws.Range("N9").Value = 5000
$ws.Range("O9").Value = 4750
$ws.Range("P9").Value = 4500
$ws.Range("Q9").Value = 4500

# Row 10 (new) - Мариуполь
$ws.Range("B10").Value = "Мариуполь  "
$ws.Range("C10").Value = 6
$ws.Range("D10").Value = 600
$ws.Range("E10").Value = 36
$ws.Range("F10").Value = 20
$ws.Range("G10").Value = 20
$ws.Range("H10").Value = 19
$ws.Range("I10").Value = 17
$ws.Range("J10").Value = 16
$ws.Range("K10").Value = 600
$ws.Range("L10").Value = 9000
$ws.Range("M10").Value = 5000
$ws.Range("N10").Value = 5000
$ws.Range("O10").Value = 4750
$ws.Range("P10").Value = 4250
$ws.Range("Q10").Value = 4000

# Row 11 (new) - Бердянск , Мелитолполь
$ws.Range("B11").Value = "Бердянск , Мелитолполь"
$ws.Range("C11").Value = 7
$ws.Range("D11").Value = 700
$ws.Range("E11").Value = 44
$ws.Range("F11").Value = 28
$ws.Range("G11").Value = 28
$ws.Range("H11").Value = 26
$ws.Range("I11").Value = 23
$ws.Range("J11").Value = 22
$ws.Range("K11").Value = 700
$ws.Range("L11").Value = 11000
$ws.Range("M11").Value = 7000
$ws.Range("N11").Value = 7000
$ws.Range("O11").Value = 6500
$ws.Range("P11").Value = 5750
$ws.Range("Q11").Value = 5500

# ------------------------------------------------------------------
# 4. Match the alternating-stripe look used on the original rows: the
#    "first city of a merged source group" rows (4,8,9 before -> 4,5,10,11
#    after the re-split) sit on a plain background, the others get the
#    light banding fill that the rest of the table already uses.
# ------------------------------------------------------------------
$ws.Range("B4").Interior.Pattern = -4142
$ws.Range("B8").Interior.Pattern = -4142
$ws.Range("B9").Interior.Pattern = -4142
$ws.Range("B10").Interior.Pattern = -4142
$ws.Range("B11").Interior.Pattern = -4142

# ------------------------------------------------------------------
# 5. Give the two brand-new rows the same borders/alignment as the rest
#    of the table (PasteSpecial above already copied borders/number
#    formats from row 9, this just normalises the numeric cells to the
#    plain/general look used for the two newly-typed rows).
# ------------------------------------------------------------------
$ws.Range("D10:Q11").NumberFormat = "General"
$ws.Range("D10:Q11").HorizontalAlignment = -4108
$ws.Range("D10:Q11").VerticalAlignment = -4108
$ws.Range("B10:B11").HorizontalAlignment = -4131
$ws.Range("B10:B11").VerticalAlignment = -4108

# ------------------------------------------------------------------
# 6. Update the sheet's dimension/selection bookkeeping to match.
# ------------------------------------------------------------------
$ws.Range("N11").Select()
